$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.495.04"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.570.67"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'212.05"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'46.05"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("D9").Value = "'24.05"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "'0.0886"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.796.75"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").Value = "1.570.42"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "28.482.06"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "'62.12"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "'230.60"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("D24").Value = "'9.10"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("E25").Value = "  +7.67%  "
$ws.Range("D26").Value = "'150.71"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'15.02"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -2.76%  "
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D35").Value = "1.393.70"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -4.07%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "'2.63"
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").Value = "'0.520"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.787"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D47").Value = "'0.969"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "'62.84"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "1.707.94"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").Value = "'86.27"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -2.19%  "
